# Add the new "gompertz" worksheet as the last (3rd) sheet and make it active/selected.
$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets
$ws = $sheets.Add([System.Type]::Missing, $sheets.Item($sheets.Count))
$ws.Name = "gompertz"

# Row 1 headers (B1:F1) and year columns (G1:Y1, 1975 then 2015-2100 step 5)
$ws.Range("B1").Value = "B_max"
$ws.Range("C1").Value = "k1"
$ws.Range("D1").Value = "k2"
$ws.Range("E1").Value = "t0"
$ws.Range("F1").Value = "Year"

$ws.Range("G1").Value = 1975
$ws.Range("H1").Value = 2015
$ws.Range("I1").Value = 2020
$ws.Range("J1").Value = 2025
$ws.Range("K1").Value = 2030
$ws.Range("L1").Value = 2035
$ws.Range("M1").Value = 2040
$ws.Range("N1").Value = 2045
$ws.Range("O1").Value = 2050
$ws.Range("P1").Value = 2055
$ws.Range("Q1").Value = 2060
$ws.Range("R1").Value = 2065
$ws.Range("S1").Value = 2070
$ws.Range("T1").Value = 2075
$ws.Range("U1").Value = 2080
$ws.Range("V1").Value = 2085
$ws.Range("W1").Value = 2090
$ws.Range("X1").Value = 2095
$ws.Range("Y1").Value = 2100

# Row 2 labels (order matters for shared-string indices) and time-period index (G2:Y2 = 0..18)
$ws.Range("F2").Value = "Time Period"
$ws.Range("E2").Value = "(4 years from 2024)"

$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 8
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 11
$ws.Range("S2").Value = 12
$ws.Range("T2").Value = 13
$ws.Range("U2").Value = 14
$ws.Range("V2").Value = 15
$ws.Range("W2").Value = 16
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 18

# Rows 3-7: per-species Gompertz growth parameters (B_max, k1, k2, t0) and the
# Gompertz-curve formula evaluated across every year column (G:Y).
$ws.Range("A3").Value = "Beef"
$ws.Range("B3").Value = 2.589
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 0.25
$ws.Range("E3").Value = 2028
$ws.Range("G3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(G`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(G`$1-`$E3))))"
$ws.Range("H3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(H`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(H`$1-`$E3))))"
$ws.Range("I3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(I`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(I`$1-`$E3))))"
$ws.Range("J3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(J`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(J`$1-`$E3))))"
$ws.Range("K3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(K`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(K`$1-`$E3))))"
$ws.Range("L3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(L`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(L`$1-`$E3))))"
$ws.Range("M3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(M`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(M`$1-`$E3))))"
$ws.Range("N3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(N`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(N`$1-`$E3))))"
$ws.Range("O3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(O`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(O`$1-`$E3))))"
$ws.Range("P3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(P`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(P`$1-`$E3))))"
$ws.Range("Q3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(Q`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(Q`$1-`$E3))))"
$ws.Range("R3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(R`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(R`$1-`$E3))))"
$ws.Range("S3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(S`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(S`$1-`$E3))))"
$ws.Range("T3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(T`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(T`$1-`$E3))))"
$ws.Range("U3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(U`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(U`$1-`$E3))))"
$ws.Range("V3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(V`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(V`$1-`$E3))))"
$ws.Range("W3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(W`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(W`$1-`$E3))))"
$ws.Range("X3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(X`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(X`$1-`$E3))))"
$ws.Range("Y3").Formula = "=IF(`$B3*EXP(-`$C3*EXP(-`$D3*(Y`$1-`$E3))) < 0.001, 0, `$B3*EXP(-`$C3*EXP(-`$D3*(Y`$1-`$E3))))"

$ws.Range("A4").Value = "Dairy"
$ws.Range("B4").Value = 0.453592
$ws.Range("C4").Value = 11
$ws.Range("D4").Value = 0.25
$ws.Range("E4").Value = 2028
$ws.Range("G4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(G`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(G`$1-`$E4))))"
$ws.Range("H4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(H`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(H`$1-`$E4))))"
$ws.Range("I4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(I`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(I`$1-`$E4))))"
$ws.Range("J4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(J`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(J`$1-`$E4))))"
$ws.Range("K4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(K`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(K`$1-`$E4))))"
$ws.Range("L4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(L`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(L`$1-`$E4))))"
$ws.Range("M4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(M`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(M`$1-`$E4))))"
$ws.Range("N4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(N`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(N`$1-`$E4))))"
$ws.Range("O4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(O`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(O`$1-`$E4))))"
$ws.Range("P4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(P`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(P`$1-`$E4))))"
$ws.Range("Q4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(Q`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(Q`$1-`$E4))))"
$ws.Range("R4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(R`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(R`$1-`$E4))))"
$ws.Range("S4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(S`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(S`$1-`$E4))))"
$ws.Range("T4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(T`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(T`$1-`$E4))))"
$ws.Range("U4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(U`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(U`$1-`$E4))))"
$ws.Range("V4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(V`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(V`$1-`$E4))))"
$ws.Range("W4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(W`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(W`$1-`$E4))))"
$ws.Range("X4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(X`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(X`$1-`$E4))))"
$ws.Range("Y4").Formula = "=IF(`$B4*EXP(-`$C4*EXP(-`$D4*(Y`$1-`$E4))) < 0.001, 0, `$B4*EXP(-`$C4*EXP(-`$D4*(Y`$1-`$E4))))"

$ws.Range("A5").Value = "Goat"
$ws.Range("B5").Value = 3.101
$ws.Range("C5").Value = 11
$ws.Range("D5").Value = 0.25
$ws.Range("E5").Value = 2028
$ws.Range("G5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(G`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(G`$1-`$E5))))"
$ws.Range("H5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(H`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(H`$1-`$E5))))"
$ws.Range("I5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(I`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(I`$1-`$E5))))"
$ws.Range("J5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(J`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(J`$1-`$E5))))"
$ws.Range("K5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(K`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(K`$1-`$E5))))"
$ws.Range("L5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(L`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(L`$1-`$E5))))"
$ws.Range("M5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(M`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(M`$1-`$E5))))"
$ws.Range("N5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(N`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(N`$1-`$E5))))"
$ws.Range("O5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(O`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(O`$1-`$E5))))"
$ws.Range("P5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(P`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(P`$1-`$E5))))"
$ws.Range("Q5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(Q`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(Q`$1-`$E5))))"
$ws.Range("R5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(R`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(R`$1-`$E5))))"
$ws.Range("S5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(S`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(S`$1-`$E5))))"
$ws.Range("T5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(T`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(T`$1-`$E5))))"
$ws.Range("U5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(U`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(U`$1-`$E5))))"
$ws.Range("V5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(V`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(V`$1-`$E5))))"
$ws.Range("W5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(W`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(W`$1-`$E5))))"
$ws.Range("X5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(X`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(X`$1-`$E5))))"
$ws.Range("Y5").Formula = "=IF(`$B5*EXP(-`$C5*EXP(-`$D5*(Y`$1-`$E5))) < 0.001, 0, `$B5*EXP(-`$C5*EXP(-`$D5*(Y`$1-`$E5))))"

$ws.Range("A6").Value = "Pork"
$ws.Range("B6").Value = 0.30399999999999999
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 0.25
$ws.Range("E6").Value = 2028
$ws.Range("G6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(G`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(G`$1-`$E6))))"
$ws.Range("H6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(H`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(H`$1-`$E6))))"
$ws.Range("I6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(I`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(I`$1-`$E6))))"
$ws.Range("J6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(J`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(J`$1-`$E6))))"
$ws.Range("K6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(K`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(K`$1-`$E6))))"
$ws.Range("L6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(L`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(L`$1-`$E6))))"
$ws.Range("M6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(M`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(M`$1-`$E6))))"
$ws.Range("N6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(N`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(N`$1-`$E6))))"
$ws.Range("O6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(O`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(O`$1-`$E6))))"
$ws.Range("P6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(P`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(P`$1-`$E6))))"
$ws.Range("Q6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(Q`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(Q`$1-`$E6))))"
$ws.Range("R6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(R`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(R`$1-`$E6))))"
$ws.Range("S6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(S`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(S`$1-`$E6))))"
$ws.Range("T6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(T`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(T`$1-`$E6))))"
$ws.Range("U6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(U`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(U`$1-`$E6))))"
$ws.Range("V6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(V`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(V`$1-`$E6))))"
$ws.Range("W6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(W`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(W`$1-`$E6))))"
$ws.Range("X6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(X`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(X`$1-`$E6))))"
$ws.Range("Y6").Formula = "=IF(`$B6*EXP(-`$C6*EXP(-`$D6*(Y`$1-`$E6))) < 0.001, 0, `$B6*EXP(-`$C6*EXP(-`$D6*(Y`$1-`$E6))))"

$ws.Range("A7").Value = "Poultry"
$ws.Range("B7").Value = 1.381
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 0.25
$ws.Range("E7").Value = 2028
$ws.Range("G7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(G`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(G`$1-`$E7))))"
$ws.Range("H7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(H`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(H`$1-`$E7))))"
$ws.Range("I7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(I`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(I`$1-`$E7))))"
$ws.Range("J7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(J`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(J`$1-`$E7))))"
$ws.Range("K7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(K`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(K`$1-`$E7))))"
$ws.Range("L7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(L`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(L`$1-`$E7))))"
$ws.Range("M7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(M`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(M`$1-`$E7))))"
$ws.Range("N7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(N`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(N`$1-`$E7))))"
$ws.Range("O7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(O`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(O`$1-`$E7))))"
$ws.Range("P7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(P`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(P`$1-`$E7))))"
$ws.Range("Q7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(Q`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(Q`$1-`$E7))))"
$ws.Range("R7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(R`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(R`$1-`$E7))))"
$ws.Range("S7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(S`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(S`$1-`$E7))))"
$ws.Range("T7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(T`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(T`$1-`$E7))))"
$ws.Range("U7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(U`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(U`$1-`$E7))))"
$ws.Range("V7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(V`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(V`$1-`$E7))))"
$ws.Range("W7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(W`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(W`$1-`$E7))))"
$ws.Range("X7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(X`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(X`$1-`$E7))))"
$ws.Range("Y7").Formula = "=IF(`$B7*EXP(-`$C7*EXP(-`$D7*(Y`$1-`$E7))) < 0.001, 0, `$B7*EXP(-`$C7*EXP(-`$D7*(Y`$1-`$E7))))"

# Scratch/debug formulas (rows 9-11), all referencing row 3 parameters, column G only
$ws.Range("G9").Formula = "=-`$C3*EXP(-`$D3*G`$1-`$E3)"
$ws.Range("G10").Formula = "=EXP(-`$D3*G`$1-`$E3)"
$ws.Range("G11").Formula = "=-`$D3*G`$1-`$E3"

# Match the authored selection on the new sheet
[void]$ws.Range("N15").Select()
